$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 18521536
$ws.Range("I106").Value = 22225044
$ws.Range("K106").Value = 22225044
$ws.Range("M106").Value = -22224413
$ws.Range("H116").Value = 2077.0667
$ws.Range("I116").Value = 2389.4
$ws.Range("J116").Value = 1452.4
$ws.Range("K116").Value = 2389.4
$ws.Range("L116").Value = 1452.4
$ws.Range("M116").Value = 1052.6
$ws.Range("N116").Value = -8336.4
$ws.Range("H125").Value = 6596184
$ws.Range("J125").Value = 11213003
$ws.Range("L125").Value = 100917027
$ws.Range("N125").Value = -100921947
$ws.Range("H138").Value = 6214440
$ws.Range("I138").Value = 2465837.5
$ws.Range("J138").Value = 7939987
$ws.Range("K138").Value = 7397512.5
$ws.Range("L138").Value = 23819961
$ws.Range("M138").Value = -7392372.5
$ws.Range("N138").Value = -23830241
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 505755
$ws.Range("I2").Value = 505755
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 505755
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -505642
$ws.Range("N2").Value = $null
$ws.Range("H32").Value = 26111.254
$ws.Range("I32").Value = 6884.25
$ws.Range("J32").Value = 117002.55
$ws.Range("K32").Value = 6884.25
$ws.Range("L32").Value = 117002.55
$ws.Range("M32").Value = -6597.25
$ws.Range("N32").Value = -117576.55
$ws.Range("H45").Value = 1527.4667
$ws.Range("I45").Value = 1624.8334
$ws.Range("J45").Value = 1138
$ws.Range("K45").Value = 1624.8334
$ws.Range("L45").Value = 1138
$ws.Range("M45").Value = -1247.8334
$ws.Range("N45").Value = -1892
$ws.Range("H88").Value = 2754.5454
$ws.Range("J88").Value = 2030
$ws.Range("L88").Value = 2030
$ws.Range("N88").Value = -2842
$ws.Range("H91").Value = 2754.5454
$ws.Range("J91").Value = 2030
$ws.Range("L91").Value = 2030
$ws.Range("N91").Value = -4838
$ws.Range("H116").Value = 505755
$ws.Range("I116").Value = 505755
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 505755
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -503461
$ws.Range("N116").Value = $null
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 253627.5
$ws.Range("I3").Value = 253627.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 253627.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -253513.5
$ws.Range("N3").Value = $null
$ws.Range("H20").Value = 2041.1666
$ws.Range("I20").Value = 1777.1111
$ws.Range("K20").Value = 1777.1111
$ws.Range("M20").Value = -1530.1111
$ws.Range("H86").Value = 1736.091
$ws.Range("I86").Value = 1947.5
$ws.Range("J86").Value = 1615.2858
$ws.Range("K86").Value = 1947.5
$ws.Range("L86").Value = 1615.2858
$ws.Range("M86").Value = -824.5
$ws.Range("N86").Value = -3861.2858
$ws.Range("H89").Value = 1736.091
$ws.Range("I89").Value = 1947.5
$ws.Range("J89").Value = 1615.2858
$ws.Range("K89").Value = 9737.5
$ws.Range("L89").Value = 8076.429
$ws.Range("M89").Value = -4121.5
$ws.Range("N89").Value = -19308.429
$ws.Range("H105").Value = 2797.4443
$ws.Range("I105").Value = 2628.1936
$ws.Range("J105").Value = 3172.2144
$ws.Range("K105").Value = 2628.1936
$ws.Range("L105").Value = 3172.2144
$ws.Range("M105").Value = -881.1936000000001
$ws.Range("N105").Value = -6666.2144
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
$ws.Range("H134").Value = 2157.4375
$ws.Range("I134").Value = 1754.8
$ws.Range("J134").Value = 4170.625
$ws.Range("K134").Value = 5264.4
$ws.Range("L134").Value = 12511.875
$ws.Range("M134").Value = -2729.4
$ws.Range("N134").Value = -17581.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 15002.5
$ws.Range("I12").Value = 16670
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 16670
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -16500
$ws.Range("N12").Value = -10340
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 922.2857
$ws.Range("J23").Value = 239.09091
$ws.Range("L23").Value = 717.27273
$ws.Range("N23").Value = -1187.27273
$ws.Range("H109").Value = 1881.2
$ws.Range("I109").Value = 1476.5
$ws.Range("K109").Value = 4429.5
$ws.Range("M109").Value = -3389.5
$ws.Range("H113").Value = 6061492.5
$ws.Range("I113").Value = 371.14285
$ws.Range("J113").Value = 11364974
$ws.Range("K113").Value = 1113.42855
$ws.Range("L113").Value = 34094922
$ws.Range("M113").Value = 1056.57145
$ws.Range("N113").Value = -34099262
$ws.Range("H131").Value = 11113019
$ws.Range("J131").Value = 12347698
$ws.Range("L131").Value = 37043094
$ws.Range("N131").Value = -37053174
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 15920
$ws.Range("J54").Value = 7400
$ws.Range("L54").Value = 7400
$ws.Range("N54").Value = -8180
$ws.Range("H80").Value = 58826310
$ws.Range("J80").Value = 166669580
$ws.Range("L80").Value = 166669580
$ws.Range("N80").Value = -166671576
$ws.Range("H83").Value = 58826310
$ws.Range("J83").Value = 166669580
$ws.Range("L83").Value = 833347900
$ws.Range("N83").Value = -833357884
$ws.Range("H132").Value = 3928.4783
$ws.Range("I132").Value = 3766.9333
$ws.Range("J132").Value = 4231.375
$ws.Range("K132").Value = 11300.7999
$ws.Range("L132").Value = 12694.125
$ws.Range("M132").Value = -8770.7999
$ws.Range("N132").Value = -17754.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2604.16
$ws.Range("I40").Value = 1760.4
$ws.Range("J40").Value = 3166.6667
$ws.Range("K40").Value = 1760.4
$ws.Range("L40").Value = 3166.6667
$ws.Range("M40").Value = -1624.4
$ws.Range("N40").Value = -3438.6667
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -26996
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -84984
$ws.Range("H125").Value = 29937.143
$ws.Range("J125").Value = 29937.143
$ws.Range("L125").Value = 29937.143
$ws.Range("N125").Value = -39777.143
$ws.Range("H131").Value = 74666.664
$ws.Range("J131").Value = 74666.664
$ws.Range("L131").Value = 74666.664
$ws.Range("N131").Value = -84746.664
$ws.Range("H132").Value = 3600.1365
$ws.Range("J132").Value = 4230.2856
$ws.Range("L132").Value = 12690.8568
$ws.Range("N132").Value = -17750.8568
